# Update countries & provincias Spain
# - Refresh the COVID numbers for 8 countries (Austria, Afganistan, Marruecos,
#   Republica de Macedonia, Eslovenia, Sri Lanka, Libano, Zambia).
# - Zambia's new total pushes it above Paraguay/Guinea Ecuatorial/Burkina Faso,
#   so the whole data table is re-sorted (descending) by "Casos totales".
# - Bump the "Datos actualizados..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-CountryRow($Country, $Casos, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes) {
    $found = $ws.Columns.Item(1).Find($Country, 1)
    $r = $found.Row()
    $ws.Cells.Item($r, 2).Value = $Casos
    $ws.Cells.Item($r, 3).Value = $Nuevos
    $ws.Cells.Item($r, 4).Value = $Activos
    $ws.Cells.Item($r, 5).Value = $Recuperados
    $ws.Cells.Item($r, 6).Value = $Criticos
    $ws.Cells.Item($r, 7).Value = $MuertesHoy
    $ws.Cells.Item($r, 8).Value = $Muertes
}

Set-CountryRow "Austria"               16353 32  14882 838  0 1 633
Set-CountryRow "Afganistan"            8145  492 930   7028 0 9 187
Set-CountryRow "Marruecos"             7048  25  4037  2817 0 1 194
Set-CountryRow "Republica de Macedonia" 1858 19  1367  381  0 4 110
Set-CountryRow "Eslovenia"             1468  1   1340  23   0 1 105
Set-CountryRow "Sri Lanka"             1027  4   584   434  0 0 9
Set-CountryRow "Libano"                961   7   251   684  0 0 26
Set-CountryRow "Zambia"                832   60  197   628  0 0 7

# Update the "last updated" footer string (last row of the table, column A)
$footer = $ws.Columns.Item(1).Find("Datos actualizados", 2)
$ws.Cells.Item($footer.Row(), 1).Value = "Datos actualizados a 20 de Mayo de 2020 a las 12:35"

# Re-sort the whole country table (A4:H219) descending by "Casos totales" (col B)
$dataRange = $ws.Range("A4:H219")
$sortKey = $ws.Range("B4:B219")
$dataRange.Sort($sortKey, 2)
